$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 1957.7
$ws.Cells.Item(86, 9).Value = 1824.2858
$ws.Cells.Item(86, 10).Value = 2269
$ws.Cells.Item(86, 11).Value = 1824.2858
$ws.Cells.Item(86, 12).Value = 2269
$ws.Cells.Item(86, 13).Value = -701.2858000000001
$ws.Cells.Item(86, 14).Value = -4515
$ws.Cells.Item(89, 8).Value = 1957.7
$ws.Cells.Item(89, 9).Value = 1824.2858
$ws.Cells.Item(89, 10).Value = 2269
$ws.Cells.Item(89, 11).Value = 9121.429
$ws.Cells.Item(89, 12).Value = 11345
$ws.Cells.Item(89, 13).Value = -3505.429
$ws.Cells.Item(89, 14).Value = -22577
$ws.Cells.Item(113, 8).Value = 27028.625
$ws.Cells.Item(113, 10).Value = 1350
$ws.Cells.Item(113, 12).Value = 1350
$ws.Cells.Item(113, 14).Value = -7858
$ws.Cells.Item(138, 8).Value = 2593.9285
$ws.Cells.Item(138, 9).Value = 2949.2856
$ws.Cells.Item(138, 10).Value = 2238.5715
$ws.Cells.Item(138, 11).Value = 8847.856800000001
$ws.Cells.Item(138, 12).Value = 6715.7145
$ws.Cells.Item(138, 13).Value = -3707.856800000001
$ws.Cells.Item(138, 14).Value = -16995.7145

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3355.6416
$ws.Cells.Item(32, 9).Value = 2013.2195
$ws.Cells.Item(32, 10).Value = 7942.25
$ws.Cells.Item(32, 11).Value = 2013.2195
$ws.Cells.Item(32, 12).Value = 7942.25
$ws.Cells.Item(32, 13).Value = -1726.2195
$ws.Cells.Item(32, 14).Value = -8516.25
$ws.Cells.Item(45, 8).Value = 2867
$ws.Cells.Item(45, 9).Value = 800
$ws.Cells.Item(45, 10).Value = 3280.4
$ws.Cells.Item(45, 11).Value = 800
$ws.Cells.Item(45, 12).Value = 3280.4
$ws.Cells.Item(45, 13).Value = -423
$ws.Cells.Item(45, 14).Value = -4034.4
$ws.Cells.Item(132, 8).Value = 2368.4827
$ws.Cells.Item(132, 9).Value = 1645.875
$ws.Cells.Item(132, 11).Value = 4937.625
$ws.Cells.Item(132, 13).Value = -2407.625

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 144984.86
$ws.Cells.Item(86, 9).Value = 2208
$ws.Cells.Item(86, 10).Value = 668500
$ws.Cells.Item(86, 11).Value = 2208
$ws.Cells.Item(86, 12).Value = 668500
$ws.Cells.Item(86, 13).Value = -1085
$ws.Cells.Item(86, 14).Value = -670746
$ws.Cells.Item(89, 8).Value = 144984.86
$ws.Cells.Item(89, 9).Value = 2208
$ws.Cells.Item(89, 10).Value = 668500
$ws.Cells.Item(89, 11).Value = 11040
$ws.Cells.Item(89, 12).Value = 3342500
$ws.Cells.Item(89, 13).Value = -5424
$ws.Cells.Item(89, 14).Value = -3353732
$ws.Cells.Item(134, 8).Value = 6696.0557
$ws.Cells.Item(134, 9).Value = 7192.074
$ws.Cells.Item(134, 10).Value = 5208
$ws.Cells.Item(134, 11).Value = 21576.222
$ws.Cells.Item(134, 12).Value = 15624
$ws.Cells.Item(134, 13).Value = -19041.222
$ws.Cells.Item(134, 14).Value = -20694

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2437.718
$ws.Cells.Item(31, 9).Value = 1005.4231
$ws.Cells.Item(31, 11).Value = 1005.4231
$ws.Cells.Item(31, 13).Value = -710.4231
$ws.Cells.Item(34, 8).Value = 2437.718
$ws.Cells.Item(34, 9).Value = 1005.4231
$ws.Cells.Item(34, 11).Value = 1005.4231
$ws.Cells.Item(34, 13).Value = -803.4231
$ws.Cells.Item(58, 8).Value = 1251.1428
$ws.Cells.Item(58, 9).Value = 1251.1428
$ws.Cells.Item(58, 11).Value = 1251.1428
$ws.Cells.Item(58, 13).Value = -1048.1428
$ws.Cells.Item(86, 8).Value = 1985.2
$ws.Cells.Item(86, 9).Value = 1985.2
$ws.Cells.Item(86, 11).Value = 1985.2
$ws.Cells.Item(86, 13).Value = -862.2
$ws.Cells.Item(89, 8).Value = 1985.2
$ws.Cells.Item(89, 9).Value = 1985.2
$ws.Cells.Item(89, 11).Value = 9926
$ws.Cells.Item(89, 13).Value = -4310
$ws.Cells.Item(132, 8).Value = 2624.8462
$ws.Cells.Item(132, 9).Value = 1323.6154
$ws.Cells.Item(132, 11).Value = 3970.8462
$ws.Cells.Item(132, 13).Value = -1440.8462
$ws.Cells.Item(134, 8).Value = 1277.6666
$ws.Cells.Item(134, 9).Value = 1285.7142
$ws.Cells.Item(134, 11).Value = 3857.1426
$ws.Cells.Item(134, 13).Value = -1322.1426
$ws.Cells.Item(136, 8).Value = 1251.1428
$ws.Cells.Item(136, 9).Value = 1251.1428
$ws.Cells.Item(136, 11).Value = 3753.4284
$ws.Cells.Item(136, 13).Value = -1203.4284

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 6421224.5
$ws.Cells.Item(131, 10).Value = 11840.5
$ws.Cells.Item(131, 12).Value = 35521.5
$ws.Cells.Item(131, 14).Value = -45601.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 38313480
$ws.Cells.Item(11, 9).Value = 7686547.5
$ws.Cells.Item(11, 10).Value = 91214540
$ws.Cells.Item(11, 11).Value = 7686547.5
$ws.Cells.Item(11, 12).Value = 91214540
$ws.Cells.Item(11, 13).Value = -7686408.5
$ws.Cells.Item(11, 14).Value = -91214818
$ws.Cells.Item(49, 8).Value = 24999
$ws.Cells.Item(49, 10).Value = 24999
$ws.Cells.Item(49, 12).Value = 24999
$ws.Cells.Item(49, 14).Value = -25367
$ws.Cells.Item(80, 8).Value = 2895.5
$ws.Cells.Item(80, 9).Value = 2895.5
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 2895.5
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).ClearContents()
$ws.Cells.Item(80, 14).Value = -1897.5
$ws.Cells.Item(83, 8).Value = 2895.5
$ws.Cells.Item(83, 9).Value = 2895.5
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 14477.5
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).ClearContents()
$ws.Cells.Item(83, 14).Value = -9485.5
$ws.Cells.Item(102, 8).Value = 1944.2307
$ws.Cells.Item(102, 9).Value = 2048.9375
$ws.Cells.Item(102, 10).Value = 1776.7
$ws.Cells.Item(102, 11).Value = 2048.9375
$ws.Cells.Item(102, 12).Value = 1776.7
$ws.Cells.Item(102, 13).Value = -426.9375
$ws.Cells.Item(102, 14).Value = -5020.7
$ws.Cells.Item(113, 8).Value = 1446.909
$ws.Cells.Item(113, 10).Value = 1564.2858
$ws.Cells.Item(113, 12).Value = 1564.2858
$ws.Cells.Item(113, 14).Value = -5904.2858
$ws.Cells.Item(132, 8).Value = 3846.4783
$ws.Cells.Item(132, 9).Value = 2418.923
$ws.Cells.Item(132, 11).Value = 7256.768999999999
$ws.Cells.Item(132, 13).Value = -4726.768999999999
$ws.Cells.Item(134, 8).Value = 44999.2
$ws.Cells.Item(134, 10).Value = 44999.2
$ws.Cells.Item(134, 12).Value = 134997.6
$ws.Cells.Item(134, 14).Value = -140067.6

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 8).Value = 14966.667
$ws.Cells.Item(20, 10).Value = 14966.667
$ws.Cells.Item(20, 12).Value = 14966.667
$ws.Cells.Item(20, 14).Value = -15418.667
$ws.Cells.Item(22, 8).Value = 1212.8334
$ws.Cells.Item(22, 9).Value = 933
$ws.Cells.Item(22, 10).Value = 1492.6666
$ws.Cells.Item(22, 11).Value = 933
$ws.Cells.Item(22, 12).Value = 1492.6666
$ws.Cells.Item(22, 13).Value = -638
$ws.Cells.Item(22, 14).Value = -2082.6666
$ws.Cells.Item(23, 8).Value = 2800
$ws.Cells.Item(23, 10).Value = 5000
$ws.Cells.Item(23, 12).Value = 5000
$ws.Cells.Item(23, 14).Value = -5460
$ws.Cells.Item(24, 8).Value = 18000
$ws.Cells.Item(24, 10).Value = 18000
$ws.Cells.Item(24, 12).Value = 18000
$ws.Cells.Item(24, 14).Value = -18686
$ws.Cells.Item(27, 8).Value = 1212.8334
$ws.Cells.Item(27, 9).Value = 933
$ws.Cells.Item(27, 10).Value = 1492.6666
$ws.Cells.Item(27, 11).Value = 933
$ws.Cells.Item(27, 12).Value = 1492.6666
$ws.Cells.Item(27, 13).Value = -826
$ws.Cells.Item(27, 14).Value = -1706.6666
$ws.Cells.Item(132, 8).Value = 1960.0869
$ws.Cells.Item(132, 9).Value = 1670.8
$ws.Cells.Item(132, 10).Value = 2040.4445
$ws.Cells.Item(132, 11).Value = 5012.4
$ws.Cells.Item(132, 12).Value = 6121.333500000001
$ws.Cells.Item(132, 13).Value = -2482.4
$ws.Cells.Item(132, 14).Value = -11181.3335

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 62609.094
$ws.Cells.Item(122, 9).Value = 81512.44
$ws.Cells.Item(122, 10).Value = 2118.4
$ws.Cells.Item(122, 11).Value = 244537.32
$ws.Cells.Item(122, 12).Value = 6355.200000000001
$ws.Cells.Item(122, 13).Value = -242087.32
$ws.Cells.Item(122, 14).Value = -11255.2
$ws.Cells.Item(126, 8).Value = 4649.4165
$ws.Cells.Item(126, 9).Value = 3393.1875
$ws.Cells.Item(126, 11).Value = 10179.5625
$ws.Cells.Item(126, 13).Value = -7709.5625
$ws.Cells.Item(132, 8).Value = 3421.111
$ws.Cells.Item(132, 9).Value = 2360.4
$ws.Cells.Item(132, 11).Value = 7081.200000000001
$ws.Cells.Item(132, 13).Value = -4551.200000000001
